# Logboek Roel Kusters - Week 7 data correction + Week 8 / Week 9 additions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the real hours for Week 7 (rows 7-14), which were all placeholder
#    zeroes before.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = 9
$ws.Range("D7").Value = 16.25
$ws.Range("F7").Value = 0.5

$ws.Range("B8").Value = 9
$ws.Range("D8").Value = 15
$ws.Range("F8").Value = 0.5

$ws.Range("B9").Value = 9
$ws.Range("D9").Value = 15
$ws.Range("F9").Value = 0.5

$ws.Range("B10").Value = 9
$ws.Range("D10").Value = 15
$ws.Range("F10").Value = 0.5

$ws.Range("B11").Value = 9
$ws.Range("D11").Value = 15
$ws.Range("F11").Value = 0.5

$ws.Range("B12").Value = 9
$ws.Range("D12").Value = 15
$ws.Range("F12").Value = 0.5

$ws.Range("B13").Value = 9
$ws.Range("D13").Value = 15
$ws.Range("F13").Value = 0.5

$ws.Range("B14").Value = 9
$ws.Range("D14").Value = 15
$ws.Range("F14").Value = 0.5

# ---------------------------------------------------------------------------
# 2. Small text corrections in the existing Week 7 task table.
# ---------------------------------------------------------------------------
$ws.Range("S15").Value = "Afstemming netwerk gedeelte van proftaak met het ander team."
$ws.Range("P24").Value = "Documentatie"
$ws.Range("S24").Value = "Logboek bijwerken."

# ---------------------------------------------------------------------------
# 3. Add the Week 7 task-table header + the new Week 8 / Week 9 task tables
#    below the existing content (rows 28-48), matching the layout used for
#    the Week 4/5/6 tables above (O=week label, P/Q/R/S = table header,
#    then one row per task).
# ---------------------------------------------------------------------------

# Copy the time-format cell style (Q/R columns) from the existing Week 6
# table so the new cells reuse the same numFmtId="20" style instead of
# Excel allocating a brand-new style entry.
$ws.Range("Q26:R26").Copy() | Out-Null

# -- Week 7 header (row 28) + rows (29-35) --------------------------------
$ws.Range("O28").Value = "Week 7"
$ws.Range("P28").Value = "Naam Taak:"
$ws.Range("Q28").Value = "Begin Tijd:"
$ws.Range("R28").Value = "Eind Tijd:"
$ws.Range("S28").Value = "Discriptie Taak:"

$ws.Range("Q29:R35").PasteSpecial(-4122) | Out-Null

$ws.Range("P29").Value = "Vergadering"
$ws.Range("Q29").Value = 0.375
$ws.Range("R29").Value = 0.39583333333333331

$ws.Range("P30").Value = "Bouwen test plc"
$ws.Range("Q30").Value = 0.39583333333333331
$ws.Range("R30").Value = 0.4375

$ws.Range("P31").Value = "Opzetten software test plc"
$ws.Range("Q31").Value = 0.4375
$ws.Range("R31").Value = 0.5

$ws.Range("P32").Value = "Pauze"
$ws.Range("Q32").Value = 0.5
$ws.Range("R32").Value = 0.52083333333333337

$ws.Range("P33").Value = "Opzetten router networking"
$ws.Range("Q33").Value = 0.52083333333333337
$ws.Range("R33").Value = 0.5625

$ws.Range("P34").Value = "Overleg met klant"
$ws.Range("Q34").Value = 0.5625
$ws.Range("R34").Value = 0.66666666666666663

$ws.Range("P35").Value = "Documentatie"
$ws.Range("Q35").Value = 0.66666666666666663
$ws.Range("R35").Value = 0.67708333333333337
$ws.Range("S35").Value = "Logboek bijwerken."

# -- Week 8 header (row 37) + rows (38-44) --------------------------------
$ws.Range("O37").Value = "Week 8"
$ws.Range("P37").Value = "Naam Taak:"
$ws.Range("Q37").Value = "Begin Tijd:"
$ws.Range("R37").Value = "Eind Tijd:"
$ws.Range("S37").Value = "Discriptie Taak:"

$ws.Range("Q38:R44").PasteSpecial(-4122) | Out-Null

$ws.Range("P38").Value = "Vergadering"
$ws.Range("Q38").Value = 0.375
$ws.Range("R38").Value = 0.40277777777777773

$ws.Range("P39").Value = "Aansturen team"
$ws.Range("Q39").Value = 0.40277777777777773
$ws.Range("R39").Value = 0.41666666666666669
$ws.Range("S39").Value = "Herverdelen taken."

$ws.Range("P40").Value = "Vergadering over server met ander team."
$ws.Range("Q40").Value = 0.41666666666666669
$ws.Range("R40").Value = 0.43055555555555558

$ws.Range("P41").Value = "Documentatie"
$ws.Range("Q41").Value = 0.43055555555555558
$ws.Range("R41").Value = 0.47916666666666669
$ws.Range("S41").Value = "Logboek bijwerken."

$ws.Range("P42").Value = "Pauze"
$ws.Range("Q42").Value = 0.47916666666666669
$ws.Range("R42").Value = 0.5

$ws.Range("P43").Value = "Onderzoek MQTT server libraries"
$ws.Range("Q43").Value = 0.5
$ws.Range("R43").Value = 0.54166666666666663

$ws.Range("P44").Value = "Onderzoek MQTT Mosquito"
$ws.Range("Q44").Value = 0.54166666666666663
$ws.Range("R44").Value = 0.625

# -- Week 9 header (row 46) + rows (47-48) --------------------------------
$ws.Range("O46").Value = "Week 9"
$ws.Range("P46").Value = "Naam Taak:"
$ws.Range("Q46").Value = "Begin Tijd:"
$ws.Range("R46").Value = "Eind Tijd:"
$ws.Range("S46").Value = "Discriptie Taak:"

$ws.Range("Q47:R48").PasteSpecial(-4122) | Out-Null

$ws.Range("P47").Value = "Vergadering"
$ws.Range("Q47").Value = 0.375
$ws.Range("R47").Value = 0.39583333333333331

$ws.Range("Q48").Value = 0.39583333333333331

# ---------------------------------------------------------------------------
# 4. Restore the selection to reflect where the author ended up after typing
#    in the new tables.
# ---------------------------------------------------------------------------
$ws.Range("P49").Select() | Out-Null
